# "Last line automatic detection"
#
# The sheet is a Key/Status/Value table. Row 13 used to be a dedicated
# "footer" row carrying a thick bottom border (marking it as the last row
# of the table). The fix removes that hard-coded footer row and instead
# lets the last data row (previously row 11: xlsx.k5 / go / xlsx.v5) sit
# two rows further down, preceded by ordinary/blank-but-styled rows - i.e.
# the table's end must be auto-detected rather than assumed to be a fixed,
# specially-bordered row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the last data row's values (Key / Status / Value) before they
# get moved two rows further down.
$lastKey    = $ws.Range("B11").Value2
$lastStatus = $ws.Range("C11").Value2
$lastValue  = $ws.Range("D11").Value2

# Row 13 carries the special "last row" thick-bottom-border formatting -
# delete it outright so that formatting disappears from the sheet.
$ws.Rows.Item(13).Delete()

# Recreate rows 13 and 14 using the plain/ordinary row format (borderId
# matching rows 5-12), by copying row 12's formatting down onto them.
$ws.Range("B12:D12").Copy()
$ws.Range("B13:D14").PasteSpecial(-4122)

# Move the last data row's values down onto the new last row (14) ...
$ws.Range("B14").Value2 = $lastKey
$ws.Range("C14").Value2 = $lastStatus
$ws.Range("D14").Value2 = $lastValue

# ... and blank out the old location (row 11), keeping its plain styling.
$ws.Range("B11:D11").ClearContents()

# The selection moves too (no longer F8).
$ws.Range("E21").Select() | Out-Null
